$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2755607.63
$ws.Range("C7").Value = -37.97976187733562
$ws.Range("D7").Value = 2802
$ws.Range("E7").Value = 2802
$ws.Range("F7").Value = 983.4431227694504
$ws.Range("G7").Value = 4.827925677708245
